$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 5857995
$ws.Cells.Item(40, 9).Value = 4918.467
$ws.Cells.Item(40, 11).Value = 4918.467
$ws.Cells.Item(40, 13).Value = -4743.467
$ws.Cells.Item(64, 8).Value = 18186044
$ws.Cells.Item(64, 9).Value = 25003936
$ws.Cells.Item(64, 10).Value = 5000
$ws.Cells.Item(64, 11).Value = 25003936
$ws.Cells.Item(64, 12).Value = 5000
$ws.Cells.Item(64, 13).Value = -25003688
$ws.Cells.Item(64, 14).Value = -5496
$ws.Cells.Item(67, 8).Value = 18186044
$ws.Cells.Item(67, 9).Value = 25003936
$ws.Cells.Item(67, 10).Value = 5000
$ws.Cells.Item(67, 11).Value = 25003936
$ws.Cells.Item(67, 12).Value = 5000
$ws.Cells.Item(67, 13).Value = -25003078
$ws.Cells.Item(67, 14).Value = -6716
$ws.Cells.Item(69, 8).Value = 16032.066
$ws.Cells.Item(69, 9).Value = 8794
$ws.Cells.Item(69, 10).Value = 16549.072
$ws.Cells.Item(69, 11).Value = 26382
$ws.Cells.Item(69, 12).Value = 49647.216
$ws.Cells.Item(69, 13).Value = -25508
$ws.Cells.Item(69, 14).Value = -51395.216
$ws.Cells.Item(72, 8).Value = 16032.066
$ws.Cells.Item(72, 9).Value = 8794
$ws.Cells.Item(72, 10).Value = 16549.072
$ws.Cells.Item(72, 11).Value = 79146
$ws.Cells.Item(72, 12).Value = 148941.648
$ws.Cells.Item(72, 13).Value = -74778
$ws.Cells.Item(72, 14).Value = -157677.648
$ws.Cells.Item(76, 8).Value = 9485.375
$ws.Cells.Item(76, 9).Value = 9346
$ws.Cells.Item(76, 11).Value = 9346
$ws.Cells.Item(76, 13).Value = -9031
$ws.Cells.Item(79, 8).Value = 9485.375
$ws.Cells.Item(79, 9).Value = 9346
$ws.Cells.Item(79, 11).Value = 9346
$ws.Cells.Item(79, 13).Value = -8254
$ws.Cells.Item(106, 8).Value = 2343.4167
$ws.Cells.Item(106, 9).Value = 2297.2856
$ws.Cells.Item(106, 11).Value = 2297.2856
$ws.Cells.Item(106, 13).Value = -1666.2856
$ws.Cells.Item(107, 8).Value = 1217.3529
$ws.Cells.Item(107, 10).Value = 1440.5
$ws.Cells.Item(107, 12).Value = 1440.5
$ws.Cells.Item(107, 14).Value = -5280.5
$ws.Cells.Item(132, 8).Value = 2839.5715
$ws.Cells.Item(132, 9).Value = 2839.5715
$ws.Cells.Item(132, 11).Value = 8518.7145
$ws.Cells.Item(132, 13).Value = -5988.7145
$ws.Cells.Item(135, 8).Value = 93751060
$ws.Cells.Item(135, 9).Value = 33334462
$ws.Cells.Item(135, 11).Value = 300010158
$ws.Cells.Item(135, 13).Value = -300007623
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 94.8
$ws.Cells.Item(5, 9).Value = 147.25
$ws.Cells.Item(5, 10).Value = 34.857143
$ws.Cells.Item(5, 11).Value = 147.25
$ws.Cells.Item(5, 12).Value = 34.857143
$ws.Cells.Item(5, 13).Value = -35.25
$ws.Cells.Item(5, 14).Value = -258.857143
$ws.Cells.Item(61, 8).Value = 32261580
$ws.Cells.Item(61, 9).Value = 35717856
$ws.Cells.Item(61, 11).Value = 35717856
$ws.Cells.Item(61, 13).Value = -35717644
$ws.Cells.Item(74, 8).Value = 27030892
$ws.Cells.Item(74, 9).Value = 35718890
$ws.Cells.Item(74, 11).Value = 35718890
$ws.Cells.Item(74, 13).Value = -35718016
$ws.Cells.Item(77, 8).Value = 27030892
$ws.Cells.Item(77, 9).Value = 35718890
$ws.Cells.Item(77, 11).Value = 178594450
$ws.Cells.Item(77, 13).Value = -178590082
$ws.Cells.Item(109, 8).Value = 80000
$ws.Cells.Item(109, 10).Value = 80000
$ws.Cells.Item(109, 12).Value = 80000
$ws.Cells.Item(109, 14).Value = -82774
$ws.Cells.Item(122, 8).Value = 2388.8484
$ws.Cells.Item(122, 9).Value = 2114.4333
$ws.Cells.Item(122, 10).Value = 5133
$ws.Cells.Item(122, 11).Value = 6343.2999
$ws.Cells.Item(122, 12).Value = 15399
$ws.Cells.Item(122, 13).Value = -3893.2999
$ws.Cells.Item(122, 14).Value = -20299
$ws.Cells.Item(132, 8).Value = 2705534.5
$ws.Cells.Item(132, 9).Value = 2859879.2
$ws.Cells.Item(132, 11).Value = 8579637.600000001
$ws.Cells.Item(132, 13).Value = -8577107.600000001
$ws.Cells.Item(136, 8).Value = 32261580
$ws.Cells.Item(136, 9).Value = 35717856
$ws.Cells.Item(136, 11).Value = 107153568
$ws.Cells.Item(136, 13).Value = -107151018
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 94.8
$ws.Cells.Item(4, 9).Value = 147.25
$ws.Cells.Item(4, 10).Value = 34.857143
$ws.Cells.Item(4, 11).Value = 147.25
$ws.Cells.Item(4, 12).Value = 34.857143
$ws.Cells.Item(4, 13).Value = -32.25
$ws.Cells.Item(4, 14).Value = -264.857143
$ws.Cells.Item(134, 8).Value = 25004344
$ws.Cells.Item(134, 10).Value = 6011.625
$ws.Cells.Item(134, 12).Value = 18034.875
$ws.Cells.Item(134, 14).Value = -23104.875
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9764.275
$ws.Cells.Item(31, 9).Value = 5633.913
$ws.Cells.Item(31, 10).Value = 15352.412
$ws.Cells.Item(31, 11).Value = 5633.913
$ws.Cells.Item(31, 12).Value = 15352.412
$ws.Cells.Item(31, 13).Value = -5338.913
$ws.Cells.Item(31, 14).Value = -15942.412
$ws.Cells.Item(34, 8).Value = 9764.275
$ws.Cells.Item(34, 9).Value = 5633.913
$ws.Cells.Item(34, 10).Value = 15352.412
$ws.Cells.Item(34, 11).Value = 5633.913
$ws.Cells.Item(34, 12).Value = 15352.412
$ws.Cells.Item(34, 13).Value = -5431.913
$ws.Cells.Item(34, 14).Value = -15756.412
$ws.Cells.Item(50, 8).Value = 60000
$ws.Cells.Item(50, 10).Value = 60000
$ws.Cells.Item(50, 12).Value = 60000
$ws.Cells.Item(50, 14).Value = -61250
$ws.Cells.Item(58, 8).Value = 50012660
$ws.Cells.Item(58, 9).Value = 83351080
$ws.Cells.Item(58, 10).Value = 5027.5
$ws.Cells.Item(58, 11).Value = 83351080
$ws.Cells.Item(58, 12).Value = 5027.5
$ws.Cells.Item(58, 13).Value = -83350877
$ws.Cells.Item(58, 14).Value = -5433.5
$ws.Cells.Item(97, 8).Value = 39996
$ws.Cells.Item(97, 10).Value = 39996
$ws.Cells.Item(97, 12).Value = 39996
$ws.Cells.Item(97, 14).Value = -41978
$ws.Cells.Item(105, 8).Value = 1819716
$ws.Cells.Item(105, 9).Value = 3334396
$ws.Cells.Item(105, 11).Value = 3334396
$ws.Cells.Item(105, 13).Value = -3332649
$ws.Cells.Item(134, 8).Value = 41668212
$ws.Cells.Item(134, 9).Value = 41668212
$ws.Cells.Item(134, 11).Value = 125004636
$ws.Cells.Item(134, 13).Value = -125002101
$ws.Cells.Item(136, 8).Value = 50012660
$ws.Cells.Item(136, 9).Value = 83351080
$ws.Cells.Item(136, 10).Value = 5027.5
$ws.Cells.Item(136, 11).Value = 250053240
$ws.Cells.Item(136, 12).Value = 15082.5
$ws.Cells.Item(136, 13).Value = -250050690
$ws.Cells.Item(136, 14).Value = -20182.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10644.5
$ws.Cells.Item(70, 9).Value = 10432.777
$ws.Cells.Item(70, 11).Value = 10432.777
$ws.Cells.Item(70, 13).Value = -10162.777
$ws.Cells.Item(73, 8).Value = 10644.5
$ws.Cells.Item(73, 9).Value = 10432.777
$ws.Cells.Item(73, 11).Value = 10432.777
$ws.Cells.Item(73, 13).Value = -9496.777
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 14).ClearContents()
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 14).ClearContents()
$ws.Cells.Item(97, 8).Value = 1380.05
$ws.Cells.Item(97, 9).Value = 843.1429000000001
$ws.Cells.Item(97, 10).Value = 2632.8333
$ws.Cells.Item(97, 11).Value = 843.1429000000001
$ws.Cells.Item(97, 12).Value = 2632.8333
$ws.Cells.Item(97, 13).Value = -347.1429000000001
$ws.Cells.Item(97, 14).Value = -3624.8333
$ws.Cells.Item(107, 8).Value = 536.7083
$ws.Cells.Item(107, 9).Value = 303.17648
$ws.Cells.Item(107, 10).Value = 1103.8572
$ws.Cells.Item(107, 11).Value = 303.17648
$ws.Cells.Item(107, 12).Value = 1103.8572
$ws.Cells.Item(107, 13).Value = 1616.82352
$ws.Cells.Item(107, 14).Value = -4943.8572
$ws.Cells.Item(126, 8).Value = 4385.8486
$ws.Cells.Item(126, 9).Value = 4385.8486
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 13157.5458
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -10687.5458
$ws.Cells.Item(126, 14).ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(95, 8).Value = 31899.5
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 13).ClearContents()
$ws.Cells.Item(96, 8).Value = 39990
$ws.Cells.Item(96, 10).Value = 39990
$ws.Cells.Item(96, 12).Value = 39990
$ws.Cells.Item(96, 14).Value = -45482
$ws.Cells.Item(99, 8).Value = 24891
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 25012950
$ws.Cells.Item(132, 9).Value = 25012950
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 75038850
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -75036320
$ws.Cells.Item(132, 14).ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 14).ClearContents()
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).ClearContents()
$ws.Cells.Item(99, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 231.89473
$ws.Cells.Item(113, 9).Value = 100.86667
$ws.Cells.Item(113, 10).Value = 723.25
$ws.Cells.Item(113, 11).Value = 302.60001
$ws.Cells.Item(113, 12).Value = 2169.75
$ws.Cells.Item(113, 13).Value = 1867.39999
$ws.Cells.Item(113, 14).Value = -6509.75
$ws.Cells.Item(126, 8).Value = 2330.9
$ws.Cells.Item(126, 9).Value = 2687.1428
$ws.Cells.Item(126, 11).Value = 8061.428400000001
$ws.Cells.Item(126, 13).Value = -5591.428400000001
$ws.Cells.Item(132, 8).Value = 15159287
$ws.Cells.Item(132, 9).Value = 19236096
$ws.Cells.Item(132, 10).Value = 16855
$ws.Cells.Item(132, 11).Value = 57708288
$ws.Cells.Item(132, 12).Value = 50565
$ws.Cells.Item(132, 13).Value = -57705758
$ws.Cells.Item(136, 8).Value = 26317412
$ws.Cells.Item(136, 9).Value = 38462404
$ws.Cells.Item(136, 10).Value = 3262.6667
$ws.Cells.Item(136, 11).Value = 115387212
$ws.Cells.Item(136, 12).Value = 9788.000100000001
$ws.Cells.Item(136, 13).Value = -115384662
$ws.Cells.Item(136, 14).Value = -14888.0001
